$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update M column (rows 3-12) with re-calculated physical distances ---
$ws.Range("M3").Value = 3.6
$ws.Range("M4").Value = 3.16
$ws.Range("M5").Value = 3
$ws.Range("M6").Value = 2.82
$ws.Range("M8").Value = 2
$ws.Range("M9").Value = 2.23
$ws.Range("M10").Value = 1.41
$ws.Range("M11").Value = 1
$ws.Range("M12").Value = 1

# --- Replace row 13 (C13:L13) with re-calculated values ---
$ws.Range("C13").Value = 3.6
$ws.Range("D13").Value = 3.16
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 2.82
$ws.Range("G13").Value = 2.23
$ws.Range("H13").Value = 2
$ws.Range("I13").Value = 2.23
$ws.Range("J13").Value = 1.41
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 1

# --- Force recalculation so the AVERAGE formula in P4 picks up new values ---
$excel.CalculateFull()

# --- Move the "#" legend entry from R9 to Q12 (clearing the old / now-unused cells) ---
$ws.Range("R9").Clear()
$ws.Range("R12").Clear()
$ws.Range("Q12").Value = "#"
$ws.Range("Q12").HorizontalAlignment = -4108
$ws.Range("Q12").Borders.LineStyle = 1

# --- Append new small distance-legend table in rows 14-17 ---
$ws.Range("O14").Value = 3.6
$ws.Range("P14").Value = 3.16
$ws.Range("Q14").Value = 3

$ws.Range("O15").Value = 2.82
$ws.Range("P15").Value = 2.23
$ws.Range("Q15").Value = 2

$ws.Range("O16").Value = 2.23
$ws.Range("P16").Value = 1.41
$ws.Range("Q16").Value = 1

$ws.Range("P17").Value = 1
$ws.Range("Q17").Value = 0

# --- Recalculate once more now that all data is in place ---
$excel.CalculateFull()

# --- Restore the workbook selection state ---
$ws.Range("K19").Select() | Out-Null
$ws.Range("K19:L19").Select() | Out-Null
